# Sync file from Google Drive
# Updates EstimatedTimeOfArrival (F), Monitored (J), OriginCode (K),
# TypeOfBus (L) and MinutesToArrival (O) values across the three
# NextBus snapshot sheets to match the refreshed source data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NextBus1")
$ws.Range("F2").Value = 45688.43523148148
$ws.Range("L2").Value = "SD"
$ws.Range("O2").Value = 17
$ws.Range("F3").Value = 45688.42366898148
$ws.Range("L3").Value = "BD"
$ws.Range("O3").Value = 1
$ws.Range("F4").Value = 45688.42335648148
$ws.Range("L4").Value = "BD"
$ws.Range("O4").Value = 0
$ws.Range("F5").Value = 45688.43210648148
$ws.Range("O5").Value = 13
$ws.Range("F6").Value = 45688.42430555556
$ws.Range("O6").Value = 2
$ws.Range("F7").Value = 45688.4255787037
$ws.Range("O7").Value = 3
$ws.Range("F8").Value = 45688.42376157407
$ws.Range("O8").Value = 1
$ws.Range("F9").Value = 45688.42681712963
$ws.Range("L9").Value = "DD"
$ws.Range("F10").Value = 45688.43042824074
$ws.Range("O10").Value = 10
$ws.Range("F11").Value = 45688.42674768518
$ws.Range("F12").Value = 45688.42503472222
$ws.Range("L12").Value = "DD"
$ws.Range("O12").Value = 3
$ws.Range("F13").Value = 45688.43229166666
$ws.Range("O13").Value = 13
$ws.Range("F14").Value = 45688.42729166667
$ws.Range("O14").Value = 6
$ws.Range("F15").Value = 45688.42604166667
$ws.Range("O15").Value = 4

$ws = $wb.Worksheets.Item("NextBus2")
$ws.Range("F2").Value = 45688.44280092593
$ws.Range("O2").Value = 28
$ws.Range("F3").Value = 45688.43153935186
$ws.Range("L3").Value = "DD"
$ws.Range("F4").Value = 45688.42685185185
$ws.Range("L4").Value = "DD"
$ws.Range("O4").Value = 5
$ws.Range("F5").Value = 45688.44219907407
$ws.Range("O5").Value = 27
$ws.Range("F6").Value = 45688.43131944445
$ws.Range("O6").Value = 12
$ws.Range("F7").Value = 45688.43327546296
$ws.Range("J7").Value = 0
$ws.Range("O7").Value = 14
$ws.Range("F8").Value = 45688.43040509259
$ws.Range("L8").Value = "DD"
$ws.Range("O8").Value = 10
$ws.Range("F9").Value = 45688.43753472222
$ws.Range("L9").Value = "SD"
$ws.Range("O9").Value = 21
$ws.Range("F10").Value = 45688.43879629629
$ws.Range("O10").Value = 22
$ws.Range("F11").Value = 45688.43148148148
$ws.Range("O11").Value = 12
$ws.Range("F12").Value = 45688.43129629629
$ws.Range("O12").Value = 12
$ws.Range("F13").Value = 45688.4419212963
$ws.Range("O13").Value = 27
$ws.Range("F14").Value = 45688.43563657408
$ws.Range("O14").Value = 18
$ws.Range("F15").Value = 45688.44493055555
$ws.Range("J15").Value = 1
$ws.Range("O15").Value = 31

$ws = $wb.Worksheets.Item("NextBus3")
$ws.Range("F2").Value = 45688.44759259259
$ws.Range("L2").Value = "DD"
$ws.Range("O2").Value = 35
$ws.Range("F3").Value = 45688.43679398148
$ws.Range("O3").Value = 20
$ws.Range("F4").Value = 45688.435
$ws.Range("L4").Value = "SD"
$ws.Range("O4").Value = 17
$ws.Range("F5").Value = 45688.44892361111
$ws.Range("J5").Value = 0
$ws.Range("O5").Value = 37
$ws.Range("F6").Value = 45688.43817129629
$ws.Range("O6").Value = 22
$ws.Range("F7").Value = 45688.44276620371
$ws.Range("O7").Value = 28
$ws.Range("F8").Value = 45688.4337037037
$ws.Range("O8").Value = 15
$ws.Range("F9").Value = 45688.44430555555
$ws.Range("L9").Value = "DD"
$ws.Range("O9").Value = 30
$ws.Range("F10").Value = 45688.4502662037
$ws.Range("O10").Value = 39
$ws.Range("F11").Value = 45688.43748842592
$ws.Range("O11").Value = 21
$ws.Range("F12").Value = 45688.43689814815
$ws.Range("O12").Value = 20
$ws.Range("F13").Value = 45688.44927083333
$ws.Range("L13").Value = "SD"
$ws.Range("F14").Value = 45688.43916666666
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 41011
$ws.Range("O14").Value = 23
$ws.Range("F15").Value = 45688.453125
$ws.Range("O15").Value = 43
